$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 3; $r -le 51; $r++) {
    $ws.Cells.Item($r, 7).Value = "Tumakuru (Tumkur)"
}

$ws.Cells.Item(24, 6).ClearContents()
